# Add the new "KI" expert response row (row 18) to the tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A18").Value = "GS_1"
$ws.Range("C18").Value = "KI"

# "Nov-11-2023" looks like a date to Excel's auto-detection, so force the
# cell to Text format before assigning it to keep it as a literal string.
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "Nov-11-2023"

$ws.Range("D18").Value = "all_submitted_trackerNov-11-2023.csv"

# Restore normal styling on the date cell so no extra formatting sticks.
$ws.Range("B18").Style = "Normal"

# Match the author's final selection recorded in the saved file.
$ws.Range("D20").Select()
